# Update generated view-count figures (column F) across sheets, per
# the "Update gh-pages to output generated at 456a3b4" regeneration.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F5").Value = 2352
$wsExhibit.Range("F16").Value = 803
$wsExhibit.Range("F20").Value = 7460
$wsExhibit.Range("F21").Value = 8393
$wsExhibit.Range("F34").Value = 1490
$wsExhibit.Range("F35").Value = 253
$wsExhibit.Range("F49").Value = 27

$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F18").Value = 311

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F8").Value = 2352
$wsAll.Range("F19").Value = 803
$wsAll.Range("F23").Value = 7460
$wsAll.Range("F24").Value = 7460
$wsAll.Range("F25").Value = 8393
$wsAll.Range("F33").Value = 253
$wsAll.Range("F50").Value = 311
$wsAll.Range("F51").Value = 27
